$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pairwise_Collab")
$ws.Rows.Item(76).Delete()
$ws.Rows.Item(72).Delete()
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(22).Delete()
